# 🔄 Actualización automática del mapa (2025-09-16 09:03:52)
# Append new claim rows (80-82) to the PEBCOM sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PEBCOM")

# Columns that must remain stored as TEXT even though they look numeric
# (Caso, F. De Reclamo, Comuna, OT) match the existing sheet convention.
$textCols = @(1,2,4,5)

$rows = @(
    @{
        A = "7224"; B = "9/16/2025"; C = "CABILDO AV. 3950"; D = "12"; E = "809784515";
        F = "PEBCOM"; G = "Pendiente"; H = "Columna inclinada"; I = 1; J = "Aplomo";
        K = "Sin equipos"; L = "Terminal"; M = -58.469735; N = -34.547232;
        O = "Saavedra"; P = "Capital Norte"
    },
    @{
        A = "7225"; B = "9/16/2025"; C = "AMENABAR 3590"; D = "13"; E = "809784519";
        F = "PEBCOM"; G = "Pendiente"; H = "Reparar rienda y tambien reclaman columna picada pero no se ve la foto verificarla y evaluar cambio"; I = 1; J = "Cambio";
        K = "Sin equipos"; L = "Terminal"; M = -58.470045; N = -34.550272;
        O = "Saavedra"; P = "Capital Norte"
    },
    @{
        A = "7234"; B = "9/16/2025"; C = "MOLDES 3388"; D = "13"; E = "809784522";
        F = "PEBCOM"; G = "Pendiente"; H = "Picada"; I = 1; J = "Cambio";
        K = "Sin equipos"; L = "Pasante"; M = -58.469426; N = -34.552639;
        O = "Saavedra"; P = "Capital Norte"
    }
)

$colOrder = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P")

$startRow = 80
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    for ($c = 0; $c -lt $colOrder.Count; $c++) {
        $colLetter = $colOrder[$c]
        $colIndex = $c + 1
        $cell = $ws.Cells.Item($r, $colIndex)
        $value = $rowData[$colLetter]
        if ($textCols -contains $colIndex) {
            # Force text storage (the source data keeps these as text,
            # not numbers), then restore the default "Normal" style so no
            # stray cell style gets attached to the new cell.
            $cell.NumberFormat = "@"
            $cell.Value = $value
            $cell.Style = "Normal"
        } else {
            $cell.Value = $value
        }
    }
}
